# Update the "Data" worksheet with new asset transfer numbers for the
# MPA test automation upload file.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Each set of 4 rows in the template repeats the same pattern of columns
# K (ANLN1), L (ANLN2), N (PANL1) that need to be bumped from the old
# asset numbers to the new ones.
$rows = @(6, 9, 11, 14, 16, 19, 21, 24, 26, 29)
foreach ($r in $rows) {
    $ws.Range("K$r").Value = 60000477
}

$rows = @(7, 8, 10, 12, 13, 15, 17, 18, 20, 22, 23, 25, 27, 28)
foreach ($r in $rows) {
    $ws.Range("L$r").Value = 327
}

$rows = @(7, 11, 12, 16, 17, 21, 22, 26, 27)
foreach ($r in $rows) {
    $ws.Range("N$r").Value = 60000478
}

$rows = @(8, 13, 18, 23, 28)
foreach ($r in $rows) {
    $ws.Range("O$r").Value = 328
}
